$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new translation-table rows (151-158) with the newly
# added English/Hebrew string pairs. (Row 159 is added further below,
# after C99, to keep the shared-string append order lined up with the
# source workbook.)
$ws.Range("B151").Value2 = "Sign In Now"
$ws.Range("C151").Value2 = "כנס עכשיו"

$ws.Range("B152").Value2 = "Sign Up Now"
$ws.Range("C152").Value2 = "הירשם עכשיו"

$ws.Range("B153").Value2 = "Don't have an account?"
$ws.Range("C153").Value2 = "אין לך עדיין חשבון?"

$ws.Range("B154").Value2 = "Sign In"
$ws.Range("C154").Value2 = "היכנס"

$ws.Range("B155").Value2 = "welcome"
$ws.Range("C155").Value2 = "ברוך הבא"

$ws.Range("B156").Value2 = "Username"
$ws.Range("C156").Value2 = "שם משתמש"

$ws.Range("B157").Value2 = "quantity"
$ws.Range("C157").Value2 = "כמות"

$ws.Range("B158").Value2 = "cost without delivery"
$ws.Range("C158").Value2 = "עלות ללא משלוח"

# --- Fix the Hebrew translation for "Deliveries" row: the old generic
# "משלוחים" string is replaced with the more specific "סוגי משלוחים"
# ("types of deliveries"), since "משלוחים" is re-purposed above for the
# "All Deliveries" row translation.
$ws.Range("C99").Value2 = "סוגי משלוחים"

$ws.Range("B159").Value2 = "Shop Now"
$ws.Range("C159").Value2 = "המשך לקנייה"

# --- Replicate the cell formatting used on the equivalent existing rows
# (column B only) by copying formats from cells that already carry the
# same style.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B153").PasteSpecial(-4122) | Out-Null

$ws.Range("B132").Copy() | Out-Null
$ws.Range("B156").PasteSpecial(-4122) | Out-Null

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B157").PasteSpecial(-4122) | Out-Null

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B159").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Match the final view state (selection) left behind in the workbook.
$ws.Range("D167").Select()
